$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# are not auto-converted to Double values (preserves exact formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "60.573.87"
$ws.Cells.Item(2, 5).Value = "  +3.44%  "

$ws.Cells.Item(3, 4).Value = "2.647.17"
$ws.Cells.Item(3, 5).Value = "  +1.19%  "

$ws.Cells.Item(4, 4).Value = "0.998"
$ws.Cells.Item(4, 5).Value = "  -0.14%  "

$ws.Cells.Item(5, 4).Value = "566.77"
$ws.Cells.Item(5, 5).Value = "  +6.13%  "

$ws.Cells.Item(6, 4).Value = "146.72"
$ws.Cells.Item(6, 5).Value = "  +2.73%  "

$ws.Cells.Item(7, 4).Value = "0.998"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "

$ws.Cells.Item(8, 4).Value = "0.611"
$ws.Cells.Item(8, 5).Value = "  +5.48%  "

$ws.Cells.Item(9, 4).Value = "2.659.81"
$ws.Cells.Item(9, 5).Value = "  +1.53%  "

$ws.Cells.Item(10, 4).Value = "6.83"
$ws.Cells.Item(10, 5).Value = "  +0.47%  "

$ws.Cells.Item(11, 5).Value = "  +4.74%  "

$ws.Cells.Item(12, 5).Value = "  +6.48%  "

$ws.Cells.Item(13, 4).Value = "0.343"
$ws.Cells.Item(13, 5).Value = "  +3.39%  "

$ws.Cells.Item(14, 4).Value = "3.109.69"
$ws.Cells.Item(14, 5).Value = "  +1.16%  "

$ws.Cells.Item(15, 4).Value = "60.538.06"
$ws.Cells.Item(15, 5).Value = "  +3.48%  "

$ws.Cells.Item(16, 4).Value = "21.91"
$ws.Cells.Item(16, 5).Value = "  +5.49%  "

$ws.Cells.Item(17, 4).Value = "0.0000137"
$ws.Cells.Item(17, 5).Value = "  +4.28%  "

$ws.Cells.Item(18, 4).Value = "2.649.97"
$ws.Cells.Item(18, 5).Value = "  +1.65%  "

$ws.Cells.Item(19, 4).Value = "4.54"
$ws.Cells.Item(19, 5).Value = "  +3.02%  "

$ws.Cells.Item(20, 4).Value = "342.62"
$ws.Cells.Item(20, 5).Value = "  +2.41%  "

$ws.Cells.Item(21, 4).Value = "10.45"
$ws.Cells.Item(21, 5).Value = "  +3.02%  "

$ws.Cells.Item(22, 4).Value = "6.36"
$ws.Cells.Item(22, 5).Value = "  +2.30%  "

$ws.Cells.Item(23, 4).Value = "5.83"
$ws.Cells.Item(23, 5).Value = "  +1.16%  "

$ws.Cells.Item(24, 5).Value = "  +0.14%  "

$ws.Cells.Item(25, 4).Value = "66.75"
$ws.Cells.Item(25, 5).Value = "  +0.11%  "

$ws.Cells.Item(26, 4).Value = "0.441"
$ws.Cells.Item(26, 5).Value = "  +4.96%  "

$ws.Cells.Item(27, 4).Value = "0.164"
$ws.Cells.Item(27, 5).Value = "  +2.96%  "

$ws.Cells.Item(28, 4).Value = "0.997"
$ws.Cells.Item(28, 5).Value = "  +0.05%  "

$ws.Cells.Item(29, 4).Value = "7.39"
$ws.Cells.Item(29, 5).Value = "  +4.27%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0805"
$ws.Cells.Item(30, 5).Value = "  +9.71%  "

$ws.Cells.Item(31, 4).Value = "0.998"
$ws.Cells.Item(31, 5).Value = "  -0.07%  "

$ws.Cells.Item(32, 4).Value = "1.71"
$ws.Cells.Item(32, 5).Value = "  +4.79%  "

$ws.Cells.Item(33, 4).Value = "6.23"
$ws.Cells.Item(33, 5).Value = "  +3.87%  "

$ws.Cells.Item(34, 4).Value = "157.84"
$ws.Cells.Item(34, 5).Value = "  +1.85%  "

$ws.Cells.Item(35, 4).Value = "19.24"
$ws.Cells.Item(35, 5).Value = "  +1.47%  "

$ws.Cells.Item(36, 4).Value = "4.10"
$ws.Cells.Item(36, 5).Value = "  +4.98%  "

$ws.Cells.Item(37, 4).Value = "0.907"
$ws.Cells.Item(37, 5).Value = "  +8.44%  "

$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).Value = "1.16"
$ws.Cells.Item(38, 5).Value = "  +5.33%  "

$ws.Cells.Item(39, 2).Value = "Fetch.AI"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(39, 4).Value = "0.898"
$ws.Cells.Item(39, 5).Value = "  +10.15%  "

$ws.Cells.Item(40, 4).Value = "37.48"
$ws.Cells.Item(40, 5).Value = "  +1.23%  "

$ws.Cells.Item(41, 5).Value = "  +6.33%  "

$ws.Cells.Item(42, 4).Value = "302.53"
$ws.Cells.Item(42, 5).Value = "  +5.70%  "

$ws.Cells.Item(43, 4).Value = "3.65"
$ws.Cells.Item(43, 5).Value = "  +1.94%  "

$ws.Cells.Item(44, 5).Value = "  -0.26%  "

$ws.Cells.Item(45, 4).Value = "0.0986"
$ws.Cells.Item(45, 5).Value = "  +4.44%  "

$ws.Cells.Item(46, 4).Value = "0.606"
$ws.Cells.Item(46, 5).Value = "  +1.63%  "

$ws.Cells.Item(47, 4).Value = "0.0547"
$ws.Cells.Item(47, 5).Value = "  +3.81%  "

$ws.Cells.Item(48, 4).Value = "128.65"
$ws.Cells.Item(48, 5).Value = "  +14.33%  "

$ws.Cells.Item(49, 4).Value = "19.38"
$ws.Cells.Item(49, 5).Value = "  +1.85%  "

$ws.Cells.Item(50, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(50, 4).Value = "10.71"
$ws.Cells.Item(50, 5).Value = "  +0.18%  "

$ws.Cells.Item(51, 2).Value = "VeChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(51, 4).Value = "0.0237"
$ws.Cells.Item(51, 5).Value = "  +5.39%  "
